$d = $word.ActiveDocument

# --- Change 1: merge "d" + "o" + "ing anything possi" runs (Persona 1 paragraph) ---
# Pure run-merge; text content is unchanged ("d"+"o"+"ing" == "doing").
$null = $d.Content.Find.Execute("firmly resolved in doing anything possi", $true, $false, $false, $false, $false, $true, 1, $false, "firmly resolved in doing anything possi", 2)

# --- Change 2: merge "Carlo and Marianna " + "go to the ... doctor " runs (Scenario 1) ---
$null = $d.Content.Find.Execute("Carlo and Marianna go to the *NOME FITTIZIO ONLUS* website to know something more about it, and to contact doctor ", $true, $false, $false, $false, $false, $true, 1, $false, "Carlo and Marianna go to the *NOME FITTIZIO ONLUS* website to know something more about it, and to contact doctor ", 2)

# --- Change 3: expand the final "." of Scenario 1 into the continuation text ---
$scenario1Para = $d.Paragraphs.Item(25)
$scenario1Range = $scenario1Para.Range
$scenario1Replacement = ". They arrive in the homepage, where they see the " + [char]0x201C + "About us" + [char]0x201D + " landmark. They click on it, and arrive to the history page of the association. After having read the history and goals of the *NOME FITTIZIO ONLUS*, they click on the " + [char]0x201C + "Services" + [char]0x201D + " landmark: from here they see the different kinds of aids the association offers, and the people involved in them. They therefore see doctor Carolinna" + [char]0x2019 + "s name under the service " + [char]0x201C + "Child neuropsychiatry aid" + [char]0x201D + ", they click on it, and arrive to the doctor" + [char]0x2019 + "s page. Here they see the doctor" + [char]0x2019 + "s professional history, and, reassured on the doctor" + [char]0x2019 + "s professional abilities, they finally reach the " + [char]0x201C + "Contact us" + [char]0x201D + " page. From here, they find the phone number of the association, and they call to fix an appointment."
$null = $scenario1Range.Find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, $scenario1Replacement, 2)

# --- Change 4: Scenario 2 paragraph text update ---
$scenario2Para = $d.Paragraphs.Item(26)
$scenario2Range = $scenario2Para.Range
$old2 = " In the meantime, Francesca phones and tells her that a play has been scheduled for this weekend at the *NOME FITTIZIO ONLUS* headquarter: so she goes to the News section to know more about it."
$new2 = " From the home page, she clicks on the landmark " + [char]0x201C + "Help us" + [char]0x201D + ". From here, she clicks on the " + [char]0x201C + "Donate now" + [char]0x201D + " link, which redirects her to an external payment gateway. After visualizing a thanks page, she gets redirected to the home page of the site. "
$null = $scenario2Range.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# --- Change 5: insert a brand-new "Scenario 3" paragraph right after Scenario 2 ---
$scenario2Para2 = $d.Paragraphs.Item(26)
$scenario2Range2 = $scenario2Para2.Range
$old3 = "After visualizing a thanks page, she gets redirected to the home page of the site. "
$scenario3Body = "Scenario 3: In the meantime of the payment, Francesca phones Roberta and tells her that a play has been scheduled for this weekend at the *NOME FITTIZIO ONLUS* headquarter. From the home page,  Roberta goes to the News section to know more about it. She arrives in the page, where she clicks on the " + [char]0x201C + "Event by date" + [char]0x201D + " link: a little calendar pops up, and she selects the dates of this week. The website redirects her to a list of events, and after scrolling she finally sees the play scheduled for the following Saturday: it is a Macbeth play. She clicks on the event title, and arrives to the event page: a beautiful picture of Lady Macbeth dominates the view, while in the description of the event, the details of date and hour are listed, along with the name of the company is going to perform. There is a ticket entrance, but the revenue of the play is all going to the association. Roberta happily phones back Francesca to tell her everything. "
$new3 = $old3 + "^p" + $scenario3Body
$null = $scenario2Range2.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# Fix up language + bold formatting on the new Scenario 3 paragraph
$scenario3Para = $d.Paragraphs.Item(27)
$scenario3Range = $scenario3Para.Range
$scenario3Range.LanguageID = "en-GB"

$scenario3BoldRange = $scenario3Para.Range
$null = $scenario3BoldRange.Find.Execute("Scenario 3: ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$scenario3BoldRange.Bold = 1

Write-Output "done"
